# Replace the C++ "EASY" row (row 3) question/choice text with the new easy-level
# C++ image filenames.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "e_c++_question.png"
$ws.Range("D3").Value = "e_c++_choiceA.png"
$ws.Range("E3").Value = "e_c++_choiceB.png"
$ws.Range("F3").Value = "e_c++_choiceC.png"
$ws.Range("G3").Value = "e_c++_choiceD.png"

# Replace the C++ "MEDIUM" row (row 4) question/choice text with the new medium-level
# C++ image filenames.
$ws.Range("C4").Value = "m_c++_question.png"
$ws.Range("D4").Value = "m_c++_choiceA.png"
$ws.Range("E4").Value = "m_c++_choiceB.png"
$ws.Range("F4").Value = "m_c++_choiceC.png"
$ws.Range("G4").Value = "m_c++_choiceD.png"

# Update the selected cell to match the author's final cursor position.
$ws.Range("H4").Select()
